$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Versión value: 0001 -> 0003
$ws.Range("B5").Value = "0003"

# Descripción value: mention the "Llamar" button
$ws.Range("B7").Value = "Se muestran en pantalla todos los datos de la oportunidad seleccionada incluyendo en la esquina inferior derecha el botón ""Llamar""."
$ws.Rows(7).RowHeight = 25.5

# Actores value: Supervisor -> Coordinador/Supervisor/Telemarketer
$ws.Range("B8").Value = "Coordinador/Supervisor/Telemarketer"

# Pre-condición value: add middle line about "Abierta" state
$ws.Range("B9").Value = "Que el actor tenga los permisos necesarios para ver el registro.`nQue existan oportunidades en estado ""Abierta""`nQue existan oportunidades asignadas a telemarketers."
$ws.Rows(9).RowHeight = 38.25

# Update the selected cell to B5
$ws.Range("B5").Select()
